# Player.xlsx fix - "fixed for some error configure file"
#
# The "Property" sheet (rows 68-75, newly added rows describing X/Y/Z/
# TargetX/TargetY/TargetZ/PathStep/LoadPropertyFinish properties) was
# missing an explicit value in the "View" (F) column, and the
# LoadPropertyFinish row (75) incorrectly had Public/Private/Save set to
# TRUE. This restores the correct boolean flags and the matching
# TRUE/FALSE list validation, then leaves the Property sheet as the
# active / selected sheet (as it was left after making the fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# --- Fix the boolean flags on rows 68-75 -----------------------------
# Rows 68-74 (X, Y, Z, TargetX, TargetY, TargetZ, PathStep): the "View"
# column (F) was left blank; it should explicitly be FALSE.
$rows = 68..74
foreach ($r in $rows) {
    $ws.Range("F$r").Value = $false
}

# Row 75 (LoadPropertyFinish): Public/Private/Save should be FALSE (not
# TRUE), and View (F) should explicitly be FALSE as well.
$ws.Range("C75").Value = $false
$ws.Range("D75").Value = $false
$ws.Range("E75").Value = $false
$ws.Range("F75").Value = $false

# --- Re-apply TRUE/FALSE list validation over the touched cells ------
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("C75:E75").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- Leave the workbook with the Property sheet active/selected ------
$ws.Activate()
$ws.Range("C75").Select()
